# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   (Overview!E/F, zh-cn!C, de-de!C all shared the same text)
# - Latest Handback DateTime refreshed for both locales
# - Error Detail cleared now that handback succeeded
# - A few columns grow/shrink to fit the new text

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

$ovw.Columns.Item(5).ColumnWidth = 29.1665
$ovw.Columns.Item(6).ColumnWidth = 29.1665

# --- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus
$zh.Range("K2").Value = "2016-09-06 14:23:36"
$zh.Range("K3").Value = "2016-09-06 14:23:36"
$zh.Range("P2").Value = ""
$zh.Range("P3").Value = ""

$zh.Columns.Item(3).ColumnWidth = 29.1665
$zh.Columns.Item(16).ColumnWidth = 12.833

# --- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus
$de.Range("K2").Value = "2016-09-06 14:23:56"
$de.Range("K3").Value = "2016-09-06 14:23:56"
$de.Range("P2").Value = ""
$de.Range("P3").Value = ""

$de.Columns.Item(3).ColumnWidth = 29.1665
$de.Columns.Item(16).ColumnWidth = 12.833
